$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 50 (the existing weekly block for this
# market/product continues below, shifted down by 3 rows).
$ws.Rows("50:52").Insert()

# Shared/fixed attributes for every row in this block.
$mercadoId = 3
$mercado = "Femacal de La Calera"
$region = "Coquimbo"
$codreg = 5
$tipo = "Fruta"
$productoId = 100107
$producto = "Otros"
$categoriaId = 100107002
$categoria = "Chirimoya"
$variedad = "Cultivar IV Región"
$unidad = "`$/bandeja 10 kilos"
$origen = "Provincia del Elquí"
$kgUnidad = 10

function Set-Fila {
    param($fila, $fecha, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $precioKg)

    $ws.Cells.Item($fila, 1).Value = $mercadoId
    $ws.Cells.Item($fila, 2).Value = $mercado
    $ws.Cells.Item($fila, 3).Value = $region
    $ws.Cells.Item($fila, 4).Value = $fecha
    $ws.Cells.Item($fila, 5).Value = $codreg
    $ws.Cells.Item($fila, 6).Value = $tipo
    $ws.Cells.Item($fila, 7).Value = $productoId
    $ws.Cells.Item($fila, 8).Value = $producto
    $ws.Cells.Item($fila, 9).Value = $categoriaId
    $ws.Cells.Item($fila, 10).Value = $categoria
    $ws.Cells.Item($fila, 11).Value = $variedad
    $ws.Cells.Item($fila, 12).Value = $calidad
    $ws.Cells.Item($fila, 13).Value = $volumen
    $ws.Cells.Item($fila, 14).Value = $precioMin
    $ws.Cells.Item($fila, 15).Value = $precioMax
    $ws.Cells.Item($fila, 16).Value = $precioProm
    $ws.Cells.Item($fila, 17).Value = $unidad
    $ws.Cells.Item($fila, 18).Value = $origen
    $ws.Cells.Item($fila, 19).Value = $precioKg
    $ws.Cells.Item($fila, 20).Value = $kgUnidad
}

Set-Fila 50 44489 "Especial" 45 30000 30000 30000 3000
Set-Fila 51 44489 "Primera"  50 25000 25000 25000 2500
Set-Fila 52 44489 "Segunda"  47 20000 20000 20000 2000
